# Add the new "Interventions maternal" worksheet, positioned right after
# "Interventions coverages" (i.e. immediately before "Interventions affected
# fraction"), and populate it with the maternal-intervention effectiveness /
# affected-fraction data. Also rename the two abbreviated intervention
# labels ("BES" / "MMS") used on the "Interventions coverages" sheet to
# their full names, matching the new wording introduced on the maternal
# sheet.

$wb = $excel.ActiveWorkbook

$coverages = $wb.Worksheets.Item("Interventions coverages")
$affected = $wb.Worksheets.Item("Interventions affected fraction")

# --- New worksheet, inserted before "Interventions affected fraction" ---
$maternal = $wb.Worksheets.Add($affected)
$maternal.Name = "Interventions maternal"

# Header row
$maternal.Range("A1").Value = "Intervention"
$maternal.Range("B1").Value = "Outcome"
$maternal.Range("C1").Value = "Pre-term SGA"
$maternal.Range("D1").Value = "Term SGA"
$maternal.Range("E1").Value = "Pre-term AGA"
$maternal.Range("F1").Value = "Term AGA"

# IPTp
$maternal.Range("A2").Value = "IPTp"
$maternal.Range("B2").Value = "effectiveness"
$maternal.Range("C2").Value = 0.35
$maternal.Range("D2").Value = 0.35
$maternal.Range("E2").Value = 0.0
$maternal.Range("F2").Value = 0.0

$maternal.Range("B3").Value = "affected fraction"
$maternal.Range("C3").Value = 0.0
$maternal.Range("D3").Value = 0.0
$maternal.Range("E3").Value = 0.0
$maternal.Range("F3").Value = 0.0

# Balanced energy supplementation (BES)
$maternal.Range("A4").Value = "Balanced energy supplementation"
$maternal.Range("B4").Value = "effectiveness"
$maternal.Range("C4").Value = 0.31
$maternal.Range("D4").Value = 0.31
$maternal.Range("E4").Value = 0.0
$maternal.Range("F4").Value = 0.0

$maternal.Range("B5").Value = "affected fraction"
$maternal.Range("C5").Value = 0.336
$maternal.Range("D5").Value = 0.336
$maternal.Range("E5").Value = 0.336
$maternal.Range("F5").Value = 0.0

# Multiple micronutrient supplementation (MMS)
$maternal.Range("A6").Value = "Multiple micronutrient supplementation"
$maternal.Range("B6").Value = "effectiveness"
$maternal.Range("C6").Value = 0.09
$maternal.Range("D6").Value = 0.09
$maternal.Range("E6").Value = 0.0
$maternal.Range("F6").Value = 0.0

$maternal.Range("B7").Value = "affected fraction"
$maternal.Range("C7").Value = 1.0
$maternal.Range("D7").Value = 1.0
$maternal.Range("E7").Value = 1.0
$maternal.Range("F7").Value = 0.0

$maternal.Range("C2:F7").NumberFormat = "#,##0.00"

# --- Rename abbreviated intervention labels on "Interventions coverages" ---
$coverages.Range("A9").Value = "Balanced energy supplementation"
$coverages.Range("A10").Value = "Multiple micronutrient supplementation"
